$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2: MAT141 material file name
$ws.Range("C2").Value = "store/materials/MAT141/gitcommands.gif"

# Update B3: course renamed from MAT111 to MILLONIARA
$ws.Range("B3").Value = "MILLONIARA"

# Update C3: material path now reflects renamed course
$ws.Range("C3").Value = "store/materials/MILLONIARA/conditional probability.pdf"

# Remove row 4 (CSC103 row) entirely, shrinking the used range to A1:C3
$ws.Rows.Item(4).Delete()
